$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$nl = [char]10

# ---------------------------------------------------------------------
# Row 3 (Tues, Jan 11) -- "Do Before Class" cell (C3).
# Swap out the "Submit substantive interest survey" bullet for a new
# "Reading Reflections Directions" bullet (the survey bullet moves down
# into the newly-populated C4 cell, with an extra note appended there).
# ---------------------------------------------------------------------
$run1 = '- Read, sign, submit syllabus on gradescope.' + $nl + '- Reading Reflections Directions <reading_reflectips.ipynb>`_' + $nl + '- Read Angrist and Piscke (MM), '
$run2 = 'Pages xi - 30' + $nl + '- `Potential Outcomes <https://github.com/nickeubank/unifyingdatascience/blob/master/lecture_slides/20_PotentialOutcomes/Fresh_Potential_Outcomes.pdf>`_'
$run3 = $nl + '- Submit Reading Reflections on Gradescope.'

$c3 = $ws.Range("C3")
$c3.Value = $run1 + $run2 + $run3

# Re-apply the original rich-text run formatting: the "Pages xi - 30..."
# chunk is bold-styled (SFBX1200) and the trailing reminder line uses the
# regular serif run font (SFRM1200), matching the untouched runs.
$c3.Characters($run1.Length + 1, $run2.Length).Font.Name = "SFBX1200"
$c3.Characters($run1.Length + $run2.Length + 1, $run3.Length).Font.Name = "SFRM1200"

# ---------------------------------------------------------------------
# Row 4 (Thurs, Jan 13) -- previously-empty "Do Before Class" cell (C4).
# Relocate the survey-submission reminder here and note reading
# reflections aren't due for this class.
# ---------------------------------------------------------------------
$ws.Range("C4").Value = '- `Submit substantive interest survey <https://forms.gle/cpr9SB4d7unXRx3j6>`_' + $nl + '- No reading reflections due.'

# The two-line entry now needs roughly double the single-line row height
# to stay fully visible (matches the row auto-growing from one line to two).
$ws.Rows(4).RowHeight = 34

# ---------------------------------------------------------------------
# Match the saved view/selection state (scrolled near the top, C7
# selected) instead of the old scrolled-to-the-bottom / C23 selection.
# ---------------------------------------------------------------------
$ws.Range("A5").Select() | Out-Null
$ws.Range("C7").Select() | Out-Null
